$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Treatment")
$ws2 = $wb.Worksheets.Item("Location")

# Copy the "Overall" summary rows (6:7) from Treatment over to Location,
# mirroring the soil-variables-only NMDS summary onto the Location sheet.
$ws1.Range("A6:P7").Copy()
$ws2.Range("A6").PasteSpecial()

# Approximate Excel's "best fit" column autosize for the newly populated
# columns on the Location sheet.
$ws2.Range("B1").ColumnWidth = 12
$ws2.Range("C1").ColumnWidth = 9
$ws2.Range("D1").ColumnWidth = 12
$ws2.Range("F1").ColumnWidth = 9
$ws2.Range("G1").ColumnWidth = 12
$ws2.Range("H1").ColumnWidth = 11
$ws2.Range("I1").ColumnWidth = 13
$ws2.Range("J1").ColumnWidth = 10
$ws2.Range("K1").ColumnWidth = 10
$ws2.Range("L1").ColumnWidth = 9
$ws2.Range("M1").ColumnWidth = 9
$ws2.Range("N1").ColumnWidth = 9
$ws2.Range("O1").ColumnWidth = 12
$ws2.Range("P1").ColumnWidth = 9

# Restore the selection left on Treatment after the copy (whole rows 6:7)
# and finish with Location active, cursor on B7 (matches the saved UI state).
$ws1.Range("A6:XFD7").Select()
$ws2.Range("B7").Select()
